# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) use numeric-looking text; force Text format on cells
# where the new value would otherwise be auto-converted to a number by Excel,
# so the stored cell content matches the source data exactly (incl. trailing zeros).
$textCells = @("D5","D6","D9","D10","D11","D12","D14","D18","D20","D21","D25","D29","D31","D32","D34","D35","D36","D37","D40","D43","D45","D46","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "57.568.39"
$ws.Range("E2").Value = "  -6.12%  "
$ws.Range("D3").Value = "2.889.65"
$ws.Range("E3").Value = "  -4.39%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "547.80"
$ws.Range("E5").Value = "  -4.10%  "
$ws.Range("D6").Value = "121.58"
$ws.Range("E6").Value = "  -5.88%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "2.880.80"
$ws.Range("E8").Value = "  -4.60%  "
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  -8.75%  "
$ws.Range("D11").Value = "4.72"
$ws.Range("E11").Value = "  -9.45%  "
$ws.Range("D12").Value = "0.432"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  -8.37%  "
$ws.Range("D14").Value = "31.58"
$ws.Range("E14").Value = "  -5.39%  "
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "3.365.53"
$ws.Range("E16").Value = "  -4.50%  "
$ws.Range("D17").Value = "2.888.62"
$ws.Range("E17").Value = "  -4.64%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "6.47"
$ws.Range("E18").Value = "  +2.56%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "57.515.14"
$ws.Range("E19").Value = "  -6.55%  "
$ws.Range("D20").Value = "406.63"
$ws.Range("E20").Value = "  -7.71%  "
$ws.Range("D21").Value = "12.79"
$ws.Range("E21").Value = "  -3.73%  "
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("E23").Value = "  -7.05%  "
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").Value = "76.75"
$ws.Range("E25").Value = "  -4.27%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -2.88%  "
$ws.Range("D29").Value = "7.17"
$ws.Range("E29").Value = "  -2.98%  "
$ws.Range("E30").Value = "  -4.14%  "
$ws.Range("D31").Value = "6.00"
$ws.Range("E31").Value = "  -5.15%  "
$ws.Range("D32").Value = "24.58"
$ws.Range("E32").Value = "  -4.21%  "
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("D34").Value = "2.02"
$ws.Range("E34").Value = "  -11.89%  "
$ws.Range("D35").Value = "0.894"
$ws.Range("E35").Value = "  -7.19%  "
$ws.Range("D36").Value = "5.32"
$ws.Range("E36").Value = "  -5.64%  "
$ws.Range("D37").Value = "48.32"
$ws.Range("E37").Value = "  -3.89%  "
$ws.Range("E38").Value = "  +7.19%  "
$ws.Range("D39").Value = "0.0₃0615"
$ws.Range("E39").Value = "  -11.35%  "
$ws.Range("D40").Value = "0.0341"
$ws.Range("E40").Value = "  -7.38%  "
$ws.Range("E41").Value = "  -4.60%  "
$ws.Range("D42").Value = "2.605.89"
$ws.Range("E42").Value = "  -2.17%  "
$ws.Range("D43").Value = "356.28"
$ws.Range("E43").Value = "  -5.08%  "
$ws.Range("E44").Value = "  -6.63%  "
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "117.39"
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("E47").Value = "  -4.14%  "
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("E49").Value = "  -2.77%  "
$ws.Range("E50").Value = "  -5.03%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.128"
$ws.Range("E51").Value = "  -3.07%  "
